# Auto-generated Excel COM-interop script
# Implements changes from commit: #10214 - set ManyToOne relations of cached entities to lazy where missing
#                                 - add cache annotation to collections of cached entities as well

$wb = $excel.ActiveWorkbook
$wsRoles = $wb.Worksheets.Item("User role")
$wsRights = $wb.Worksheets.Item("User Rights")

# --- Update Uuid column (User role sheet, column FT) ---
$wsRoles.Range("FT2").Value = 'QJCP7Z-AD2A2D-YJWYXL-7BLW2D3I'
$wsRoles.Range("FT3").Value = 'RSNKSP-LUEDRG-BZEH5O-KEZO2EFI'
$wsRoles.Range("FT4").Value = 'TBIGFY-OFUONA-UTN5OR-EQ6FCBJM'
$wsRoles.Range("FT5").Value = 'TNHILE-XKAZ5L-NSIDFH-TZ3UKIYY'
$wsRoles.Range("FT6").Value = 'UTQNR6-VRTKWZ-U42SEQ-DDJZCKBI'
$wsRoles.Range("FT7").Value = 'V52HST-R5GAPK-KUS2FB-IY2OKHFI'
$wsRoles.Range("FT8").Value = 'WDNLR4-QODXES-TA3EDM-42CM2P2Y'
$wsRoles.Range("FT9").Value = 'S7ANPX-YN3Q2R-MEPIQG-IWGQCOFU'
$wsRoles.Range("FT10").Value = 'XTYY5L-YU3OJF-YYVEGP-AHEYKDNU'
$wsRoles.Range("FT11").Value = 'X7S5GU-HF3FBA-6CJ7XA-LKBD2E2A'
$wsRoles.Range("FT12").Value = 'VTKRYY-XKJCZV-6UWMIJ-45ZVSA74'
$wsRoles.Range("FT13").Value = 'WBZTYE-N4RVK4-XZMUFZ-NFXOSETQ'
$wsRoles.Range("FT14").Value = 'SINIIY-7EKXCP-A6R3BX-WEV52A2U'
$wsRoles.Range("FT15").Value = 'QNWFWI-4CHEFD-WTOGKD-7LP3CG7M'
$wsRoles.Range("FT16").Value = 'RVLNRH-MU5HSI-AJZMXS-BXRB2GQI'
$wsRoles.Range("FT17").Value = 'TA37XI-I5ITUM-ORVBNZ-76DUCHDQ'
$wsRoles.Range("FT18").Value = 'QDO6BR-DUDGQ6-7CH5B6-IQ3WSIMA'
$wsRoles.Range("FT19").Value = 'XUKH5K-HL6NQH-QXCG2N-M42GKCHI'
$wsRoles.Range("FT20").Value = 'XBEFLD-IZKOTV-XB6DHT-3ZM4KINI'
$wsRoles.Range("FT21").Value = 'QPC4DB-ZPE6RR-CARHKB-KBKJ2MV4'
$wsRoles.Range("FT22").Value = 'TBSJAA-LICJEI-UK5INB-2LFECNXE'
$wsRoles.Range("FT23").Value = 'UH74FA-L4N3SM-4ZISGJ-RTXN2FJ4'
$wsRoles.Range("FT24").Value = 'T3SETM-J6GC4Q-VPKEBV-TEMSSFJ4'
$wsRoles.Range("FT25").Value = 'WRPTSW-SMWVNR-4EOG3X-JY7W2JU4'
$wsRoles.Range("FT26").Value = 'U3INAZ-FOAO5T-TDFBVH-MFVAKFLU'
$wsRoles.Range("FT27").Value = 'QYFL24-R2JSCQ-PHLA4I-MQL2CLSE'

# --- Update Needed user rights column (User Rights sheet, column E) ---
$wsRights.Range("E4").Value = 'PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$wsRights.Range("E6").Value = 'ADDITIONAL_TEST_VIEW, DOCUMENT_VIEW, THERAPY_VIEW, TASK_VIEW, SAMPLE_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, CASE_VIEW, TREATMENT_DELETE, PRESCRIPTION_DELETE, TASK_DELETE, CLINICAL_VISIT_DELETE, IMMUNIZATION_VIEW, PERSON_VIEW, PATHOGEN_TEST_DELETE, CLINICAL_COURSE_VIEW, DOCUMENT_DELETE, VISIT_DELETE, SAMPLE_VIEW, IMMUNIZATION_DELETE'
$wsRights.Range("E9").Value = 'CASE_EDIT, PERSON_VIEW, PERSON_EDIT, CASE_VIEW'
$wsRights.Range("E19").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW'
$wsRights.Range("E20").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW, PERSON_EDIT'
$wsRights.Range("E22").Value = 'IMMUNIZATION_VIEW, PERSON_VIEW, VISIT_DELETE, PERSON_DELETE'
$wsRights.Range("E27").Value = 'PERSON_VIEW, PERSON_EDIT'
$wsRights.Range("E32").Value = 'ADDITIONAL_TEST_VIEW, PATHOGEN_TEST_DELETE, SAMPLE_VIEW, ADDITIONAL_TEST_DELETE'
$wsRights.Range("E34").Value = 'SAMPLE_EDIT, SAMPLE_VIEW'
$wsRights.Range("E41").Value = 'ADDITIONAL_TEST_VIEW, SAMPLE_VIEW'
$wsRights.Range("E45").Value = 'PERSON_VIEW, CONTACT_VIEW, CASE_VIEW'
$wsRights.Range("E46").Value = 'PERSON_VIEW, PERSON_EDIT, CONTACT_VIEW, CASE_VIEW'
$wsRights.Range("E48").Value = 'ADDITIONAL_TEST_VIEW, DOCUMENT_VIEW, TASK_VIEW, SAMPLE_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, CASE_VIEW, TASK_DELETE, PERSON_VIEW, PATHOGEN_TEST_DELETE, DOCUMENT_DELETE, VISIT_DELETE, SAMPLE_VIEW, CONTACT_VIEW'
$wsRights.Range("E51").Value = 'CASE_CREATE, PERSON_VIEW, PERSON_EDIT, CONTACT_EDIT, CASE_VIEW, CONTACT_VIEW'
$wsRights.Range("E52").Value = 'PERSON_VIEW, PERSON_EDIT, CONTACT_EDIT, CONTACT_VIEW, CASE_VIEW'
$wsRights.Range("E66").Value = 'DOCUMENT_VIEW, EVENT_VIEW, DOCUMENT_DELETE'
$wsRights.Range("E72").Value = 'ADDITIONAL_TEST_VIEW, DOCUMENT_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW, TASK_VIEW, EVENTPARTICIPANT_DELETE, ACTION_DELETE, SAMPLE_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE, TASK_DELETE, PERSON_VIEW, PATHOGEN_TEST_DELETE, DOCUMENT_DELETE, SAMPLE_VIEW, VISIT_DELETE'
$wsRights.Range("E75").Value = 'EVENT_VIEW, EVENT_EDIT'
$wsRights.Range("E77").Value = 'PERSON_VIEW, EVENT_VIEW'
$wsRights.Range("E78").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW'
$wsRights.Range("E79").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW, PERSON_EDIT'
$wsRights.Range("E81").Value = 'PERSON_VIEW, ADDITIONAL_TEST_VIEW, PATHOGEN_TEST_DELETE, EVENTPARTICIPANT_VIEW, EVENT_VIEW, SAMPLE_VIEW, VISIT_DELETE, SAMPLE_DELETE, PERSON_DELETE, ADDITIONAL_TEST_DELETE'
$wsRights.Range("E83").Value = 'PERSON_VIEW, EVENTPARTICIPANT_VIEW, EVENT_VIEW, PERSON_EDIT, EVENTPARTICIPANT_EDIT'
$wsRights.Range("E106").Value = 'PERSON_VIEW, DASHBOARD_CONTACT_VIEW, CONTACT_VIEW, CASE_VIEW'
$wsRights.Range("E110").Value = 'PERSON_VIEW, THERAPY_VIEW, CASE_VIEW'
$wsRights.Range("E117").Value = 'PERSON_VIEW, CLINICAL_COURSE_VIEW, THERAPY_VIEW, CASE_VIEW'
$wsRights.Range("E122").Value = 'PERSON_VIEW, PORT_HEALTH_INFO_VIEW, CASE_VIEW'
$wsRights.Range("E142").Value = 'PERSON_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS'
$wsRights.Range("E143").Value = 'PERSON_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, TRAVEL_ENTRY_VIEW'
$wsRights.Range("E144").Value = 'PERSON_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, TRAVEL_ENTRY_VIEW, PERSON_EDIT'
$wsRights.Range("E146").Value = 'TASK_DELETE, PERSON_VIEW, DOCUMENT_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, DOCUMENT_DELETE, TASK_VIEW, TRAVEL_ENTRY_VIEW, VISIT_DELETE, PERSON_DELETE'
$wsRights.Range("E165").Value = 'CASE_EDIT, EVENTPARTICIPANT_CREATE, EVENTPARTICIPANT_VIEW, EVENT_CREATE, EVENTPARTICIPANT_EDIT, PERSON_DELETE, EXTERNAL_MESSAGE_VIEW, PATHOGEN_TEST_CREATE, IMMUNIZATION_VIEW, IMMUNIZATION_EDIT, PERSON_VIEW, SAMPLE_EDIT, VISIT_DELETE, EVENT_EDIT, SAMPLE_CREATE, CONTACT_VIEW, CASE_CREATE, EVENT_VIEW, CONTACT_EDIT, PATHOGEN_TEST_EDIT, CASE_VIEW, PATHOGEN_TEST_DELETE, CONTACT_CREATE, PERSON_EDIT, SAMPLE_VIEW, IMMUNIZATION_CREATE, IMMUNIZATION_DELETE'
